$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 1000000
$ws.Range("I9").Value = 1000000
$ws.Range("M9").Value = -999831
$ws.Range("K9").Value = 1000000

# Row 11
$ws.Range("M11").Value = -737.8461
$ws.Range("K11").Value = 877.8461
$ws.Range("I11").Value = 877.8461
$ws.Range("H11").Value = 877.8461

# Row 18
$ws.Range("K18").Value = 1790.5555
$ws.Range("H18").Value = 1790.5555
$ws.Range("M18").Value = -1506.5555
$ws.Range("I18").Value = 1790.5555

# Row 62
$ws.Range("K62").Value = 20837206
$ws.Range("I62").Value = 20837206
$ws.Range("H62").Value = 19234636
$ws.Range("M62").Value = -20836582

# Row 65
$ws.Range("H65").Value = 19234636
$ws.Range("K65").Value = 104186030
$ws.Range("I65").Value = 20837206
$ws.Range("M65").Value = -104182910

# Row 69
$ws.Range("K69").Value = 60000
$ws.Range("M69").Value = -59126
$ws.Range("H69").Value = 86208
$ws.Range("I69").Value = 20000

# Row 72
$ws.Range("H72").Value = 86208
$ws.Range("M72").Value = -175632
$ws.Range("K72").Value = 180000
$ws.Range("I72").Value = 20000

# Row 107
$ws.Range("M107").Value = 1280.4
$ws.Range("K107").Value = 639.6
$ws.Range("H107").Value = 647.9091
$ws.Range("I107").Value = 639.6

# Row 116
$ws.Range("J116").Value = 14886.429
$ws.Range("N116").Value = -21770.429
$ws.Range("L116").Value = 14886.429
$ws.Range("H116").Value = 16554.5

# Row 121
$ws.Range("L121").Value = 2997
$ws.Range("J121").Value = 999
$ws.Range("H121").Value = 999
$ws.Range("N121").Value = -6491

# Row 132
$ws.Range("I132").Value = 3485.9143
$ws.Range("H132").Value = 3880.762
$ws.Range("M132").Value = -7927.742899999999
$ws.Range("K132").Value = 10457.7429

# Row 134
$ws.Range("J134").Value = 57999
$ws.Range("H134").Value = 57999
$ws.Range("L134").Value = 57999
$ws.Range("N134").Value = -68139

# Row 137
$ws.Range("K137").Value = 3657.2223
$ws.Range("M137").Value = -1107.2223
$ws.Range("I137").Value = 1219.0741
$ws.Range("H137").Value = 1460.4333

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("I32").Value = 932.8049
$ws.Range("K32").Value = 932.8049
$ws.Range("M32").Value = -645.8049
$ws.Range("H32").Value = 1307.9546

# Row 76
$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20676

# Row 79
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("H79").Value = 20000
$ws.Range("N79").Value = -22340

# Row 86
$ws.Range("I86").Value = 44000
$ws.Range("K86").Value = 44000
$ws.Range("H86").Value = 44000
$ws.Range("M86").Value = -42814

# Row 89
$ws.Range("I89").Value = 44000
$ws.Range("H89").Value = 44000
$ws.Range("K89").Value = 132000
$ws.Range("M89").Value = -126072

$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 2624.7646
$ws.Range("L64").Value = 765.375
$ws.Range("N64").Value = -1215.375
$ws.Range("J64").Value = 765.375

# Row 67
$ws.Range("H67").Value = 2624.7646
$ws.Range("L67").Value = 765.375
$ws.Range("N67").Value = -2325.375
$ws.Range("J67").Value = 765.375

# Row 100
$ws.Range("H100").Value = 9747.5
$ws.Range("J100").Value = 9747.5
$ws.Range("L100").Value = 9747.5
$ws.Range("N100").Value = -11911.5

# Row 105
$ws.Range("M105").Value = -879.5
$ws.Range("K105").Value = 2626.5
$ws.Range("I105").Value = 2626.5
$ws.Range("H105").Value = 2751.2856

# Row 107
$ws.Range("M107").Value = 803.75
$ws.Range("K107").Value = 1116.25
$ws.Range("H107").Value = 1416
$ws.Range("I107").Value = 1116.25

# Row 132
$ws.Range("J132").Value = 74075.60000000001
$ws.Range("N132").Value = -84195.60000000001
$ws.Range("H132").Value = 74075.60000000001
$ws.Range("L132").Value = 74075.60000000001

# Row 133
$ws.Range("L133").Value = 95000
$ws.Range("H133").Value = 95000
$ws.Range("J133").Value = 95000
$ws.Range("N133").Value = -105120

# Row 134
$ws.Range("M134").Value = -25808.823
$ws.Range("K134").Value = 28343.823
$ws.Range("H134").Value = 9184.9375
$ws.Range("I134").Value = 9447.941000000001

# Row 139
$ws.Range("H139").Value = 250000
$ws.Range("L139").Value = 250000
$ws.Range("J139").Value = 250000
$ws.Range("N139").Value = -260280

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("I7").Value = 139.06667
$ws.Range("L7").Value = 132.71428
$ws.Range("H7").Value = 137.04546
$ws.Range("M7").Value = -26.06666999999999
$ws.Range("K7").Value = 139.06667
$ws.Range("J7").Value = 132.71428
$ws.Range("N7").Value = -358.71428

# Row 22
$ws.Range("I22").Value = 395.33334
$ws.Range("K22").Value = 395.33334
$ws.Range("L22").Value = 1493.2
$ws.Range("N22").Value = -2193.2
$ws.Range("H22").Value = 1310.2222
$ws.Range("M22").Value = -45.33334000000002
$ws.Range("J22").Value = 1493.2

# Row 31
$ws.Range("I31").Value = 1782
$ws.Range("K31").Value = 1782
$ws.Range("H31").Value = 3211.28
$ws.Range("M31").Value = -1487

# Row 34
$ws.Range("K34").Value = 1782
$ws.Range("M34").Value = -1580
$ws.Range("H34").Value = 3211.28
$ws.Range("I34").Value = 1782

# Row 58
$ws.Range("K58").Value = 4147.273
$ws.Range("M58").Value = -3944.273
$ws.Range("H58").Value = 5163.222
$ws.Range("I58").Value = 4147.273

# Row 99
$ws.Range("K99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("M99").ClearContents()

# Row 126
$ws.Range("I126").Value = 0
$ws.Range("H126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

# Row 134
$ws.Range("M134").Value = -13212.8775
$ws.Range("K134").Value = 15747.8775
$ws.Range("H134").Value = 5950.413
$ws.Range("I134").Value = 5249.2925

# Row 136
$ws.Range("M136").Value = -9891.819
$ws.Range("I136").Value = 4147.273
$ws.Range("K136").Value = 12441.819
$ws.Range("H136").Value = 5163.222

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("K97").Value = 607.46155
$ws.Range("M97").Value = -111.46155
$ws.Range("I97").Value = 607.46155
$ws.Range("L97").Value = 200561.2
$ws.Range("J97").Value = 200561.2
$ws.Range("H97").Value = 56150.168
$ws.Range("N97").Value = -201553.2

# Row 135
$ws.Range("H135").Value = 93999
$ws.Range("L135").Value = 93999
$ws.Range("N135").Value = -104139
$ws.Range("J135").Value = 93999

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 23343.375
$ws.Range("N122").Value = -28243.375
$ws.Range("J122").Value = 7781.125
$ws.Range("H122").Value = 7781.125
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 33
$ws.Range("H33").Value = 11245
$ws.Range("L33").Value = 13490
$ws.Range("J33").Value = 13490
$ws.Range("N33").Value = -13990

# Row 36
$ws.Range("H36").Value = 11245
$ws.Range("L36").Value = 13490
$ws.Range("J36").Value = 13490
$ws.Range("N36").Value = -13990

# Row 37
$ws.Range("H37").Value = 39500
$ws.Range("K37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("M37").ClearContents()

# Row 68
$ws.Range("N68").Value = -27906.143
$ws.Range("H68").Value = 26284.143
$ws.Range("L68").Value = 26284.143
$ws.Range("J68").Value = 26284.143

# Row 69
$ws.Range("J69").Value = 33125
$ws.Range("H69").Value = 33125
$ws.Range("L69").Value = 33125
$ws.Range("N69").Value = -34623

# Row 71
$ws.Range("J71").Value = 26284.143
$ws.Range("N71").Value = -86964.429
$ws.Range("H71").Value = 26284.143
$ws.Range("L71").Value = 78852.429

# Row 72
$ws.Range("H72").Value = 33125
$ws.Range("N72").Value = -106863
$ws.Range("L72").Value = 99375
$ws.Range("J72").Value = 33125

# Row 97
$ws.Range("L97").Value = 20594.889
$ws.Range("J97").Value = 20594.889
$ws.Range("H97").Value = 20594.889
$ws.Range("N97").Value = -22576.889

# Row 126
$ws.Range("I126").Value = 3051.1428
$ws.Range("M126").Value = -6683.428400000001
$ws.Range("H126").Value = 3429.25
$ws.Range("K126").Value = 9153.428400000001
